$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Testmail #12: Ik heb nog geen geld terug."
$ws.Range("B8").Value = "Beste klant,`nBedankt voor uw bericht. Om uw terugbetaling te kunnen verwerken, heb ik wat meer informatie nodig. Kunt u mij uw ordernummer of transactiereferentie geven, zodat ik het voor u kan nakijken?`nAlvast bedankt voor uw medewerking.`nMet vriendelijke groet,`n[Naam] `nE-mailassistent"
$ws.Range("C8").Value = "Ik heb nog geen geld terug."
$ws.Range("D8").Value = "mailmind.test@zohomail.eu"
$ws.Range("E8").Value = "Retour / Terugbetaling"
$ws.Range("F8").Value = "2025-08-04 20:51:30"
$ws.Range("G8").Value = "Ja"
$ws.Range("H8").Value = "Nee"
$ws.Range("I8").Value = "Ja"
$ws.Range("J8").Value = "Nee"

$ws.Rows.Item(8).AutoFit()
